$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.404.89"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "1.820.13"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.36"
$ws.Range("D5").ClearFormats()

$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5122"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3920"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07811"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.75"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.108"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.95"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.240"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.56%  "

$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.481"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.14%  "

$ws.Range("D16").Value = "1.822.05"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001133"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.49"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06619"
$ws.Range("D19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.69"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.082"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").Value = "28.439.89"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.261"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.09"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.60%  "

$ws.Range("D27").Value = "2.029.71"
$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.394"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.94%  "

$ws.Range("E30").Value = "  +1.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1100"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.102"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.669"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.646"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07051"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2212"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02321"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.184"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.775"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6251"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.21"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.176"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9996"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("E44").Value = "  -0.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.39"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("E46").Value = "  +0.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5877"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.43"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.973"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.193"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06890"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.02%  "
